# Fixed blackout date and last facility inspection date.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Previous Inspection Date" (G10) and "Most Recent Inspection
# Date" (H10) for the MORAN VISTA ASSISTED LIVING row. Dependent formulas
# in I10, J10, L10, M10 recalc automatically.
$ws.Range("G10").Value = [DateTime]"2016-07-24"
$ws.Range("H10").Value = [DateTime]"2018-03-12"

# Update the saved selection to match what was active when the file was
# last saved.
$ws.Activate()
$ws.Range("E14").Select()
